$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- 1. Insert three new bullet paragraphs before the
#        "Developed and deployed custom analytical tools..." bullet ---
$rng = $d.Content
$rng.Find.Execute("Developed and deployed custom analytical tools and algorithms", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    $insertPoint = $rng.Paragraphs(1).Range
    $insertPoint.Collapse(1)

    $newText1 = "$bullet Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data"
    $newText2 = "$bullet Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters"
    $newText3 = "$bullet Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts"

    $insertPoint.InsertBefore($newText1 + [char]13 + $newText2 + [char]13 + $newText3 + [char]13)
}

# --- 2. Remove the "Created fraud detection systems for campaign finance..." bullet ---
$rng2 = $d.Content
$rng2.Find.Execute("Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng2.Find.Found) {
    $para = $rng2.Paragraphs(1)
    $para.Range.Delete()
}
